$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.004.87'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '1.859.65'
$ws.Range('E3').Value = '  -0.48%  '
$c = $ws.Range('D4')
$c.Value = "'1.003"
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.26%  '
$c = $ws.Range('D5')
$c.Value = "'311.88"
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.Value = "'0.5082"
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('E8').Value = '  +0.07%  '
$c = $ws.Range('D9')
$c.Value = "'0.08255"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -7.54%  '
$ws.Range('E10').Value = '  -0.82%  '
$c = $ws.Range('D11')
$c.Value = "'41.51"
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +0.14%  '
$c = $ws.Range('D12')
$c.Value = "'6.197"
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -1.73%  '
$c = $ws.Range('D13')
$c.Value = "'20.47"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('D14').Value = '1.858.33'
$ws.Range('E14').Value = '  -0.18%  '
$c = $ws.Range('D15')
$c.Value = "'7.183"
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.68%  '
$c = $ws.Range('D16')
$c.Value = "'1.002"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.15%  '
$c = $ws.Range('D17')
$c.Value = "'0.00001095"
$c.Style = 'Normal'
$ws.Range('E17').Value = '  -0.38%  '
$c = $ws.Range('D18')
$c.Value = "'90.44"
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.42%  '
$c = $ws.Range('D19')
$c.Value = "'0.06604"
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.45%  '
$c = $ws.Range('D20')
$c.Value = "'17.72"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('E21').Value = '  +0.20%  '
$c = $ws.Range('D22')
$c.Value = "'6.013"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('D23').Value = '28.029.56'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('E24').Value = '  -3.78%  '
$c = $ws.Range('D25')
$c.Value = "'2.239"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.94%  '
$c = $ws.Range('D26')
$c.Value = "'2.548"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +1.41%  '
$ws.Range('D27').Value = '2.072.33'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('E29').Value = '  -1.48%  '
$c = $ws.Range('D30')
$c.Value = "'124.38"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('E32').Value = '  -1.90%  '
$c = $ws.Range('D33')
$c.Value = "'5.591"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +0.33%  '
$c = $ws.Range('D34')
$c.Value = "'3.598"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +0.27%  '
$c = $ws.Range('D35')
$c.Value = "'9.591"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +2.82%  '
$c = $ws.Range('D36')
$c.Value = "'0.06528"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.22%  '
$c = $ws.Range('D37')
$c.Value = "'0.02408"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.56%  '
$ws.Range('E38').Value = '  -1.07%  '
$c = $ws.Range('D39')
$c.Value = "'1.204"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.51%  '
$ws.Range('E40').Value = '  -3.93%  '
$c = $ws.Range('D41')
$c.Value = "'0.6393"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('E42').Value = '  -3.74%  '
$c = $ws.Range('D43')
$c.Value = "'4.864"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.72%  '
$c = $ws.Range('D44')
$c.Value = "'0.6072"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.27%  '
$c = $ws.Range('D45')
$c.Value = "'13.03"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('E46').Value = '  -0.12%  '
$c = $ws.Range('D47')
$c.Value = "'3.662"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.12%  '
$c = $ws.Range('D48')
$c.Value = "'1.978"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('E49').Value = '  -1.78%  '
$c = $ws.Range('D50')
$c.Value = "'120.63"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.74%  '
$c = $ws.Range('D51')
$c.Value = "'78.90"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +1.13%  '
